# Auto-generated cell updates for scheduled Sheets refresh (Bahamut_Profits)
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2450.8333
$ws.Range("I98").Value = 2536.7273
$ws.Range("K98").Value = 2536.7273
$ws.Range("M98").Value = -1038.7273
$ws.Range("H121").Value = 724.94446
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 715.5625
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 2146.6875
$ws.Range("M121").Value = -653
$ws.Range("N121").Value = -5640.6875
$ws.Range("H122").Value = 2450.8333
$ws.Range("I122").Value = 2536.7273
$ws.Range("K122").Value = 7610.1819
$ws.Range("M122").Value = -5160.1819
$ws.Range("H129").Value = 931.04346
$ws.Range("J129").Value = 1108.8572
$ws.Range("L129").Value = 3326.5716
$ws.Range("N129").Value = -13326.5716
$ws.Range("H132").Value = 3228444.5
$ws.Range("I132").Value = 3070.7144
$ws.Range("J132").Value = 10001730
$ws.Range("K132").Value = 9212.143199999999
$ws.Range("L132").Value = 30005190
$ws.Range("M132").Value = -6682.143199999999
$ws.Range("N132").Value = -30010250
$ws.Range("H141").Value = 1561.5625
$ws.Range("I141").Value = 1561.5625
$ws.Range("K141").Value = 4684.6875
$ws.Range("M141").Value = 495.3125

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3455.47
$ws.Range("I32").Value = 3455.47
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3455.47
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3168.47
$ws.Range("H61").Value = 1582.5714
$ws.Range("I61").Value = 1618.129
$ws.Range("J61").Value = 1307
$ws.Range("K61").Value = 1618.129
$ws.Range("L61").Value = 1307
$ws.Range("M61").Value = -1406.129
$ws.Range("N61").Value = -1731
$ws.Range("H63").Value = 2211.3928
$ws.Range("I63").Value = 2135.8262
$ws.Range("J63").Value = 2559
$ws.Range("K63").Value = 2135.8262
$ws.Range("L63").Value = 2559
$ws.Range("M63").Value = -1449.8262
$ws.Range("N63").Value = -3931
$ws.Range("H66").Value = 2211.3928
$ws.Range("I66").Value = 2135.8262
$ws.Range("J66").Value = 2559
$ws.Range("K66").Value = 10679.131
$ws.Range("L66").Value = 12795
$ws.Range("M66").Value = -7247.130999999999
$ws.Range("N66").Value = -19659
$ws.Range("H74").Value = 2200.4443
$ws.Range("I74").Value = 2263
$ws.Range("J74").Value = 1700
$ws.Range("K74").Value = 2263
$ws.Range("L74").Value = 1700
$ws.Range("M74").Value = -1389
$ws.Range("N74").Value = -3448
$ws.Range("H77").Value = 2200.4443
$ws.Range("I77").Value = 2263
$ws.Range("J77").Value = 1700
$ws.Range("K77").Value = 11315
$ws.Range("L77").Value = 8500
$ws.Range("M77").Value = -6947
$ws.Range("N77").Value = -17236
$ws.Range("H122").Value = 1006
$ws.Range("I122").Value = 1006
$ws.Range("K122").Value = 3018
$ws.Range("M122").Value = -568
$ws.Range("H126").Value = 8200
$ws.Range("I126").Value = 8200
$ws.Range("K126").Value = 24600
$ws.Range("M126").Value = -22130
$ws.Range("H136").Value = 1582.5714
$ws.Range("I136").Value = 1618.129
$ws.Range("J136").Value = 1307
$ws.Range("K136").Value = 4854.387
$ws.Range("L136").Value = 3921
$ws.Range("M136").Value = -2304.387
$ws.Range("N136").Value = -9021

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 29698.695
$ws.Range("I20").Value = 44284.066
$ws.Range("J20").Value = 2351.125
$ws.Range("K20").Value = 44284.066
$ws.Range("L20").Value = 2351.125
$ws.Range("M20").Value = -44037.066
$ws.Range("N20").Value = -2845.125
$ws.Range("H86").Value = 45457572
$ws.Range("I86").Value = 71431144
$ws.Range("J86").Value = 3824
$ws.Range("K86").Value = 71431144
$ws.Range("L86").Value = 3824
$ws.Range("M86").Value = -71430021
$ws.Range("N86").Value = -6070
$ws.Range("H89").Value = 45457572
$ws.Range("I89").Value = 71431144
$ws.Range("J89").Value = 3824
$ws.Range("K89").Value = 357155720
$ws.Range("L89").Value = 19120
$ws.Range("M89").Value = -357150104
$ws.Range("N89").Value = -30352
$ws.Range("H99").Value = 200003100
$ws.Range("I99").Value = 200003100
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 200003100
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -200001602

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37601.395
$ws.Range("I31").Value = 1516.95
$ws.Range("J31").Value = 127812.5
$ws.Range("K31").Value = 1516.95
$ws.Range("L31").Value = 127812.5
$ws.Range("M31").Value = -1221.95
$ws.Range("N31").Value = -128402.5
$ws.Range("H34").Value = 37601.395
$ws.Range("I34").Value = 1516.95
$ws.Range("J34").Value = 127812.5
$ws.Range("K34").Value = 1516.95
$ws.Range("L34").Value = 127812.5
$ws.Range("M34").Value = -1314.95
$ws.Range("N34").Value = -128216.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 337
$ws.Range("I38").Value = 81.75
$ws.Range("J38").Value = 628.7143
$ws.Range("K38").Value = 245.25
$ws.Range("L38").Value = 1886.1429
$ws.Range("M38").Value = 101.75
$ws.Range("N38").Value = -2580.1429
$ws.Range("H58").Value = 3769.95
$ws.Range("J58").Value = 3957.842
$ws.Range("L58").Value = 11873.526
$ws.Range("N58").Value = -12129.526
$ws.Range("H121").Value = 881.93335
$ws.Range("I121").Value = 100
$ws.Range("J121").Value = 937.7857
$ws.Range("K121").Value = 300
$ws.Range("L121").Value = 2813.3571
$ws.Range("M121").Value = 1010
$ws.Range("N121").Value = -5433.3571
$ws.Range("H131").Value = 19308966
$ws.Range("I131").Value = 55666930
$ws.Range("J131").Value = 60632
$ws.Range("K131").Value = 167000790
$ws.Range("L131").Value = 181896
$ws.Range("M131").Value = -166995750
$ws.Range("N131").Value = -191976

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4059
$ws.Range("I70").Value = 3653.8462
$ws.Range("K70").Value = 3653.8462
$ws.Range("M70").Value = -3383.8462
$ws.Range("H73").Value = 4059
$ws.Range("I73").Value = 3653.8462
$ws.Range("K73").Value = 3653.8462
$ws.Range("M73").Value = -2717.8462
$ws.Range("H113").Value = 1680.875
$ws.Range("I113").Value = 995.6667
$ws.Range("J113").Value = 2092
$ws.Range("K113").Value = 995.6667
$ws.Range("L113").Value = 2092
$ws.Range("M113").Value = 1174.3333
$ws.Range("N113").Value = -6432
$ws.Range("H122").Value = 1298.3334
$ws.Range("I122").Value = 972.5
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 2917.5
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -467.5
$ws.Range("N122").Value = -10750
$ws.Range("H126").Value = 2792.9167
$ws.Range("I126").Value = 2950.1
$ws.Range("J126").Value = 2007
$ws.Range("K126").Value = 8850.299999999999
$ws.Range("L126").Value = 6021
$ws.Range("M126").Value = -6380.299999999999
$ws.Range("N126").Value = -10961

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 339.08334
$ws.Range("I55").Value = 371.2857
$ws.Range("J55").Value = 294
$ws.Range("K55").Value = 371.2857
$ws.Range("L55").Value = 294
$ws.Range("M55").Value = -198.2857
$ws.Range("N55").Value = -640
$ws.Range("H122").Value = 34629.332
$ws.Range("I122").Value = 34629.332
$ws.Range("K122").Value = 103887.996
$ws.Range("M122").Value = -101437.996

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 55557280
$ws.Range("I122").Value = 62501676
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 187505028
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -187502578
$ws.Range("N122").Value = -11200
